$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "juliandate"
$ws.Range("A13").Value = "datetime "
$ws.Range("D14").Value = "In Progress"
$ws.Range("G12").Value = "proper syntax in .txt file "
$ws.Range("G13").Value = "proper syntax in .txt file "

$ws.Range("C15").Select()
